$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Cells.Item(2, 5).Value = 3
    $ws.Cells.Item(2, 7).Value = 339.4154793333333
    $ws.Cells.Item(2, 8).Value = 1018.246438
    $ws.Cells.Item(2, 9).Value = 0.5189044318468032
    $ws.Cells.Item(2, 10).Value = 0.5189044318468032
    $ws.Cells.Item(2, 11).Value = 3
    $ws.Cells.Item(2, 13).Value = 21.08181366666667
    $ws.Cells.Item(2, 14).Value = 63.245441
    $ws.Cells.Item(2, 15).Value = 0.0571606014598545
    $ws.Cells.Item(2, 16).Value = 0.0571606014598545
    $ws.Cells.Item(2, 17).Value = 7155.493890887684
    $ws.Cells.Item(2, 18).Value = 64399.44501798916
    $ws.Cells.Item(2, 19).Value = 0.02966088942454735
    $ws.Cells.Item(2, 20).Value = 0.02966088942454735
    $ws.Cells.Item(3, 5).Value = 3
    $ws.Cells.Item(3, 7).Value = 339.4154793333333
    $ws.Cells.Item(3, 8).Value = 1018.246438
    $ws.Cells.Item(3, 9).Value = 0.5189044318468032
    $ws.Cells.Item(3, 10).Value = 0.5189044318468032
    $ws.Cells.Item(3, 11).Value = 3
    $ws.Cells.Item(3, 13).Value = 301.6001486666667
    $ws.Cells.Item(3, 14).Value = 904.800446
    $ws.Cells.Item(3, 15).Value = 0.8177496571571792
    $ws.Cells.Item(3, 16).Value = 0.8177496571571792
    $ws.Cells.Item(3, 17).Value = 102367.7590267013
    $ws.Cells.Item(3, 18).Value = 921309.8312403114
    $ws.Cells.Item(3, 19).Value = 0.4243339212400642
    $ws.Cells.Item(3, 20).Value = 0.4243339212400642
    $ws.Cells.Item(4, 5).Value = 3
    $ws.Cells.Item(4, 7).Value = 339.4154793333333
    $ws.Cells.Item(4, 8).Value = 1018.246438
    $ws.Cells.Item(4, 9).Value = 0.5189044318468032
    $ws.Cells.Item(4, 10).Value = 0.5189044318468032
    $ws.Cells.Item(4, 11).Value = 3
    $ws.Cells.Item(4, 13).Value = 46.13524966666667
    $ws.Cells.Item(4, 14).Value = 138.405749
    $ws.Cells.Item(4, 15).Value = 0.1250897413829664
    $ws.Cells.Item(4, 16).Value = 0.1250897413829664
    $ws.Cells.Item(4, 17).Value = 15659.01787977468
    $ws.Cells.Item(4, 18).Value = 140931.1609179721
    $ws.Cells.Item(4, 19).Value = 0.06490962118219171
    $ws.Cells.Item(4, 20).Value = 0.06490962118219171
    $ws.Cells.Item(5, 5).Value = 3
    $ws.Cells.Item(5, 7).Value = 243.634776
    $ws.Cells.Item(5, 8).Value = 730.9043280000001
    $ws.Cells.Item(5, 9).Value = 0.3724731861573274
    $ws.Cells.Item(5, 10).Value = 0.3724731861573274
    $ws.Cells.Item(5, 11).Value = 3
    $ws.Cells.Item(5, 13).Value = 21.08181366666667
    $ws.Cells.Item(5, 14).Value = 63.245441
    $ws.Cells.Item(5, 15).Value = 0.0571606014598545
    $ws.Cells.Item(5, 16).Value = 0.0571606014598545
    $ws.Cells.Item(5, 17).Value = 5136.262950352072
    $ws.Cells.Item(5, 18).Value = 46226.36655316865
    $ws.Cells.Item(5, 19).Value = 0.02129079134842119
    $ws.Cells.Item(5, 20).Value = 0.02129079134842119
    $ws.Cells.Item(6, 5).Value = 3
    $ws.Cells.Item(6, 7).Value = 243.634776
    $ws.Cells.Item(6, 8).Value = 730.9043280000001
    $ws.Cells.Item(6, 9).Value = 0.3724731861573274
    $ws.Cells.Item(6, 10).Value = 0.3724731861573274
    $ws.Cells.Item(6, 11).Value = 3
    $ws.Cells.Item(6, 13).Value = 301.6001486666667
    $ws.Cells.Item(6, 14).Value = 904.800446
    $ws.Cells.Item(6, 15).Value = 0.8177496571571792
    $ws.Cells.Item(6, 16).Value = 0.8177496571571792
    $ws.Cells.Item(6, 17).Value = 73480.28466197003
    $ws.Cells.Item(6, 18).Value = 661322.5619577303
    $ws.Cells.Item(6, 19).Value = 0.3045898202803967
    $ws.Cells.Item(6, 20).Value = 0.3045898202803967
    $ws.Cells.Item(7, 5).Value = 3
    $ws.Cells.Item(7, 7).Value = 243.634776
    $ws.Cells.Item(7, 8).Value = 730.9043280000001
    $ws.Cells.Item(7, 9).Value = 0.3724731861573274
    $ws.Cells.Item(7, 10).Value = 0.3724731861573274
    $ws.Cells.Item(7, 11).Value = 3
    $ws.Cells.Item(7, 13).Value = 46.13524966666667
    $ws.Cells.Item(7, 14).Value = 138.405749
    $ws.Cells.Item(7, 15).Value = 0.1250897413829664
    $ws.Cells.Item(7, 16).Value = 0.1250897413829664
    $ws.Cells.Item(7, 17).Value = 11240.15121824241
    $ws.Cells.Item(7, 18).Value = 101161.3609641817
    $ws.Cells.Item(7, 19).Value = 0.04659257452850957
    $ws.Cells.Item(7, 20).Value = 0.04659257452850957
    $ws.Cells.Item(8, 5).Value = 3
    $ws.Cells.Item(8, 7).Value = 71.049919
    $ws.Cells.Item(8, 8).Value = 213.149757
    $ws.Cells.Item(8, 9).Value = 0.1086223819958692
    $ws.Cells.Item(8, 10).Value = 0.1086223819958692
    $ws.Cells.Item(8, 11).Value = 3
    $ws.Cells.Item(8, 13).Value = 21.08181366666667
    $ws.Cells.Item(8, 14).Value = 63.245441
    $ws.Cells.Item(8, 15).Value = 0.0571606014598545
    $ws.Cells.Item(8, 16).Value = 0.0571606014598545
    $ws.Cells.Item(8, 17).Value = 1497.86115338976
    $ws.Cells.Item(8, 18).Value = 13480.75038050784
    $ws.Cells.Item(8, 19).Value = 0.006208920686885957
    $ws.Cells.Item(8, 20).Value = 0.006208920686885957
    $ws.Cells.Item(9, 5).Value = 3
    $ws.Cells.Item(9, 7).Value = 71.049919
    $ws.Cells.Item(9, 8).Value = 213.149757
    $ws.Cells.Item(9, 9).Value = 0.1086223819958692
    $ws.Cells.Item(9, 10).Value = 0.1086223819958692
    $ws.Cells.Item(9, 11).Value = 3
    $ws.Cells.Item(9, 13).Value = 301.6001486666667
    $ws.Cells.Item(9, 14).Value = 904.800446
    $ws.Cells.Item(9, 15).Value = 0.8177496571571792
    $ws.Cells.Item(9, 16).Value = 0.8177496571571792
    $ws.Cells.Item(9, 17).Value = 21428.66613315462
    $ws.Cells.Item(9, 18).Value = 192857.9951983916
    $ws.Cells.Item(9, 19).Value = 0.08882591563671823
    $ws.Cells.Item(9, 20).Value = 0.08882591563671823
    $ws.Cells.Item(10, 5).Value = 3
    $ws.Cells.Item(10, 7).Value = 71.049919
    $ws.Cells.Item(10, 8).Value = 213.149757
    $ws.Cells.Item(10, 9).Value = 0.1086223819958692
    $ws.Cells.Item(10, 10).Value = 0.1086223819958692
    $ws.Cells.Item(10, 11).Value = 3
    $ws.Cells.Item(10, 13).Value = 46.13524966666667
    $ws.Cells.Item(10, 14).Value = 138.405749
    $ws.Cells.Item(10, 15).Value = 0.1250897413829664
    $ws.Cells.Item(10, 16).Value = 0.1250897413829664
    $ws.Cells.Item(10, 17).Value = 3277.905751861445
    $ws.Cells.Item(10, 18).Value = 29501.151766753
    $ws.Cells.Item(10, 19).Value = 0.01358754567226506
    $ws.Cells.Item(10, 20).Value = 0.01358754567226506
